# Applies the "for delivery to add data" edit:
#  1. Updates the "Document Created on ..." date.
#  2. Clears every bolded answer value in the Gen-Info list (the
#     labels like "Medical Record Number:  " stay, only the bold
#     run holding the answer is emptied out, keeping the empty
#     run + its <w:rPr><w:b/></w:rPr> so formatting survives).

$d = $word.ActiveDocument

# 1) Update the creation date line.
$d.Content.Find.Execute(
    "Document Created on 14-Mar-2018", $true, $false, $false, $false, $false,
    $true, 1, $false, "Document Created on 02-Apr-2018", 2) | Out-Null

# 2) Clear the bold "answer" run in each of these ListBullet paragraphs.
#    (Paragraph indices are stable across the edit - we only ever empty
#    an existing run's text, never add/remove paragraphs.)
$targetIndices = @(
    3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,19,20,21,22,24,25,26,28,29,30,
    31,32,33,34,35,36,37,38,39,40,41,42,43,45,46,47,48,49,50,51,52,53,54,
    55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,
    93,94,96,97,98,99,100,101,102,103,104
)

# Minimal WordOpenXML package wrapping a single empty bold run - used to
# overwrite the answer range so the bold run survives with no text,
# matching "<w:r><w:rPr><w:b/></w:rPr></w:r>".
$emptyBoldRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:b/></w:rPr></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

foreach ($idx in $targetIndices) {
    $p = $d.Paragraphs.Item($idx)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End - 1   # exclude the paragraph mark

    # Find where the bold answer run begins by scanning forward for the
    # first bold character (the label run itself is never bold).
    $valueStart = -1
    for ($i = $pStart; $i -lt $pEnd; $i++) {
        $ch = $d.Range($i, $i + 1)
        if ($ch.Font.Bold -eq -1) {
            $valueStart = $i
            break
        }
    }

    if ($valueStart -ge 0 -and $valueStart -lt $pEnd) {
        $valueRange = $d.Range($valueStart, $pEnd)
        $valueRange.InsertXML($emptyBoldRunXml)
    }
}
